$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update the rolled-up metrics after trade #20 closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.93               # Current Capital
$summary.Range("B4").Value = -0.07000000000000001  # Total P&L $
$summary.Range("B5").Value = -0.07000000000000001  # Total P&L %
$summary.Range("B6").Value = 20                     # Total Trades
$summary.Range("B7").Value = 6                      # Winning Trades
$summary.Range("B9").Value = 30                     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 4).
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.93000000000001   # Capital
$status.Range("D4").Value = 20                  # Trades
$status.Range("E4").Value = -0.07000000000000001 # P&L $
$status.Range("F4").Value = -0.07000000000000001 # P&L %
$status.Range("G4").Value = 30                  # Win Rate %

# ---------------------------------------------------------------------------
# Append trade #20 to both the "All Trades" sheet and the "MarketMaking"
# sheet (row 21 on each).
# ---------------------------------------------------------------------------
$newRow = @(20, "2026-02-17", "07:59:45", "MarketMaking", "DOWN", 0.97, 0.98, "CLOSED", 1.0309, 0.01, 99.93000000000001, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.14)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Item(21, 1).Value = $newRow[0]
    # B21 holds a literal date-formatted string ("2026-02-17"), not a real
    # date. Plain assignment lets Excel auto-coerce it into a date serial,
    # so force a text format, assign, then clear the format override back
    # off (the source file keeps this cell on the default/general style).
    $ws.Cells.Item(21, 2).NumberFormat = "@"
    $ws.Cells.Item(21, 2).Value = $newRow[1]
    $ws.Cells.Item(21, 2).ClearFormats()
    $ws.Cells.Item(21, 3).Value = $newRow[2]
    $ws.Cells.Item(21, 4).Value = $newRow[3]
    $ws.Cells.Item(21, 5).Value = $newRow[4]
    $ws.Cells.Item(21, 6).Value = $newRow[5]
    $ws.Cells.Item(21, 7).Value = $newRow[6]
    $ws.Cells.Item(21, 8).Value = $newRow[7]
    $ws.Cells.Item(21, 9).Value = $newRow[8]
    $ws.Cells.Item(21, 10).Value = $newRow[9]
    $ws.Cells.Item(21, 11).Value = $newRow[10]
    $ws.Cells.Item(21, 12).Value = $newRow[11]
    $ws.Cells.Item(21, 13).Value = $newRow[12]
    $ws.Cells.Item(21, 14).Value = $newRow[13]
    $ws.Cells.Item(21, 15).Value = $newRow[14]
    $ws.Cells.Item(21, 16).Value = $newRow[15]
    $ws.Cells.Item(21, 17).Value = $newRow[16]
}
